$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the value to be stored as plain text (matches the shared-string
    # "t=s" cells already used throughout this workbook for things that look
    # like numbers/dates, e.g. "64.00" or "2019-09-04"), then drop the
    # temporary text number-format again so the cell keeps its original
    # (unstyled) look.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1) "Unexpected" sheet: the "new trait for marie__1" columns (E:F) were
#    recorded here by mistake; they belong on "Stress_management" instead.
#    Deleting them shifts "2:Soil_porosity_%__2" / "1:Soil_porosity_%__1"
#    (and their TIMESTAMP columns) left from G:J into E:H.
# ---------------------------------------------------------------------------
$wsUnexpected = $wb.Worksheets.Item("Unexpected")
$wsUnexpected.Range("E1:F5").Delete(-4159)  # xlShiftToLeft

# ---------------------------------------------------------------------------
# 2) "Stress_management" sheet: add the "new trait for marie__1" columns
#    (G:H) that were removed from "Unexpected" above, with the same data.
# ---------------------------------------------------------------------------
$wsStress = $wb.Worksheets.Item("Stress_management")

# Copy the existing header formatting (bold font + border) from E1 onto the
# two new header cells so they match the rest of the header row exactly.
$wsStress.Range("E1").Copy() | Out-Null
$wsStress.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsStress.Range("G1").Value = "new trait for marie__1"
$wsStress.Range("H1").Value = "TIMESTAMP_new trait for marie__1"

$newTraitValues = @{ 2 = "hello"; 3 = "salut"; 4 = "Ola"; 5 = "hallo" }
foreach ($r in 2..5) {
    Set-TextValue $wsStress.Range("G$r") $newTraitValues[$r]
    Set-TextValue $wsStress.Range("H$r") "2019-09-04"
}

# ---------------------------------------------------------------------------
# 3) "for_analysis" sheet: columns AC:AJ mirror the two sheets above, so
#    reorder them the same way: 1:Soil_porosity, Biotic_stress_control,
#    new trait for marie, 2:Soil_porosity.
# ---------------------------------------------------------------------------
$wsAnalysis = $wb.Worksheets.Item("for_analysis")

$wsAnalysis.Range("AC1").Value = "1:Soil_porosity_%__1"
$wsAnalysis.Range("AD1").Value = "TIMESTAMP_1:Soil_porosity_%__1"
$wsAnalysis.Range("AE1").Value = "Biotic_stress_control_product_amount_g/m2__1"
$wsAnalysis.Range("AF1").Value = "TIMESTAMP_Biotic_stress_control_product_amount_g/m2__1"
$wsAnalysis.Range("AG1").Value = "new trait for marie__1"
$wsAnalysis.Range("AH1").Value = "TIMESTAMP_new trait for marie__1"
$wsAnalysis.Range("AI1").Value = "2:Soil_porosity_%__2"
$wsAnalysis.Range("AJ1").Value = "TIMESTAMP_2:Soil_porosity_%__2"

$analysisData = @{
    2 = @("64.00", "2019-08-29", "10.00", "2019-09-04", "hello", "2019-09-04", "63.00", "2019-08-29")
    3 = @("55.00", "2019-08-29", "2.00",  "2019-09-04", "salut", "2019-09-04", "52.00", "2019-08-29")
    4 = @("63.00", "2019-08-29", "88.00", "2019-09-04", "Ola",   "2019-09-04", "54.00", "2019-08-29")
    5 = @("63.00", "2019-08-29", "7.00",  "2019-09-04", "hallo", "2019-09-04", "62.00", "2019-08-29")
}
$analysisCols = @("AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ")
foreach ($r in 2..5) {
    for ($i = 0; $i -lt $analysisCols.Length; $i++) {
        Set-TextValue $wsAnalysis.Range("$($analysisCols[$i])$r") $analysisData[$r][$i]
    }
}
